# Update "Förändrad" date (column C) for rows 2-8 from 2023-10-25 (45224)
# to 2023-11-03 (45233), leaving formatting/styles untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
